$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write a literal text value into a cell without Excel's "looks like a
# number" auto-conversion kicking in (and without leaving a stray
# NumberFormat style behind): build it as a quoted-string formula, then
# paste-special just the resulting value back over itself.
function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# The log keeps the latest 6 equations. A new equation was evaluated, so
# the oldest row (5) is dropped, rows 6-7 shift up into 5-6, and the new
# result is appended at row 7.
Set-TextValue $ws.Cells.Item(5,1) $ws.Cells.Item(6,1).Text
Set-TextValue $ws.Cells.Item(5,2) $ws.Cells.Item(6,2).Text
Set-TextValue $ws.Cells.Item(5,3) $ws.Cells.Item(6,3).Text

Set-TextValue $ws.Cells.Item(6,1) $ws.Cells.Item(7,1).Text
Set-TextValue $ws.Cells.Item(6,2) $ws.Cells.Item(7,2).Text
Set-TextValue $ws.Cells.Item(6,3) $ws.Cells.Item(7,3).Text

Set-TextValue $ws.Cells.Item(7,1) "2+2"
Set-TextValue $ws.Cells.Item(7,2) "4"
Set-TextValue $ws.Cells.Item(7,3) "1648571968260"

$excel.CutCopyMode = $false

# Header row no longer carries an explicit per-cell style override.
$ws.Range("A1:C1").Style = "Normal"

# Columns B and C were widened (no longer auto-fit to content).
$ws.Columns.Item(2).ColumnWidth = 21.0221354166667
$ws.Columns.Item(3).ColumnWidth = 20.3072916666667

$ws.Range("G4").Select()
